$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was added to the series. Insert a fresh
# row at row 47 (pushing the existing rows 47-140 down to 48-141) and
# populate it with the new observation.
$ws.Rows.Item(47).Insert()

$ws.Range("A47").Value = 4
$ws.Range("B47").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C47").Value = "Los Lagos"
$ws.Range("D47").Value = 44838
$ws.Range("E47").Value = 10
$ws.Range("F47").Value = 100112022
$ws.Range("G47").Value = "Arveja Verde"
$ws.Range("H47").Value = "Perfection"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 70
$ws.Range("K47").Value = 35000
$ws.Range("L47").Value = 35000
$ws.Range("M47").Value = 35000
$ws.Range("N47").Value = "`$/malla 25 kilos"
$ws.Range("O47").Value = "Provincia de Huasco"
$ws.Range("P47").Value = 1400
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

# Match the source D-column number format (date serial display).
$ws.Range("D47").NumberFormat = $ws.Range("D48").NumberFormat
